$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The font used by the totals in column T (rows 2-7) grows from 10pt to 11pt.
$ws.Range("T2:T7").Font.Size = 11

# Updated totals in column T.
$ws.Range("T2").Value = 141
$ws.Range("T3").Value = 37
$ws.Range("T4").Value = 859
$ws.Range("T5").Value = 46
$ws.Range("T6").Value = 7

# Move the active selection to O9 (was U5).
$ws.Range("O9").Select()
